$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update asset JSON text values (bug fix of ids) ---

# Arun's assets (row 2, col G)
$ws.Range("G2").Value = "[{""id"":""YB-99 "",""type"":""Monitor""},{""id"":""YB-92-L"",""type"":""Laptop""},{""id"":""YB-92"",""type"":""Mouse""},{""id"":""YB-39-H"",""type"":""Headphone""}]`n"

# hp's assets (row 3, col G)
$ws.Range("G3").Value = "[{""id"":""YB-19 "",""type"":""Monitor""},{""id"":""YB-19"",""type"":""Laptop""},{""id"":""YB-19"",""type"":""Mouse""},{""id"":""YB-19"",""type"":""Headphone""}]`n"

# Hariharan's assets (row 4, col G)
$ws.Range("G4").Value = "[{""id"":""YB-92 "",""type"":""Monitor""},{""id"":""YB-11-L"",""type"":""Laptop""},{""id"":""YB-11-MS"",""type"":""Mouse""},{""id"":""YB-11-H"",""type"":""Headphone""}]`n"

# --- Fix formatting: wrap text on G3 and grow row 3 height to fit ---
$ws.Range("G3").WrapText = $true
$ws.Range("G3").Font.Name = "Calibri"
$ws.Range("G3").Font.ThemeColor = 1
$ws.Rows.Item(3).RowHeight = 34.5
